# Auto-generated update of market price columns (H-N) across multiple sheets.
# Mirrors a scheduled-runner refresh of currentAveragePrice / LevePrice / LeveProfit
# columns for specific Leve rows, as captured by the upstream OOXML diff.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 60.272728
$ws.Range("I39").Value = 51.75
$ws.Range("J39").Value = 83
$ws.Range("K39").Value = 155.25
$ws.Range("L39").Value = 249
$ws.Range("M39").Value = 140.75
$ws.Range("N39").Value = -841
$ws.Range("H86").Value = 62504776
$ws.Range("I86").Value = 4156.222
$ws.Range("J86").Value = 142862720
$ws.Range("K86").Value = 4156.222
$ws.Range("L86").Value = 142862720
$ws.Range("M86").Value = -3033.222
$ws.Range("N86").Value = -142864966
$ws.Range("H89").Value = 62504776
$ws.Range("I89").Value = 4156.222
$ws.Range("J89").Value = 142862720
$ws.Range("K89").Value = 20781.11
$ws.Range("L89").Value = 714313600
$ws.Range("M89").Value = -15165.11
$ws.Range("N89").Value = -714324832
$ws.Range("H116").Value = 4227.864
$ws.Range("I116").Value = 4301.353
$ws.Range("J116").Value = 3978
$ws.Range("K116").Value = 4301.353
$ws.Range("L116").Value = 3978
$ws.Range("M116").Value = -859.3530000000001
$ws.Range("N116").Value = -10862
$ws.Range("H121").Value = 487.5263
$ws.Range("J121").Value = 470.16666
$ws.Range("L121").Value = 1410.49998
$ws.Range("N121").Value = -4904.499980000001
$ws.Range("H129").Value = 1059286.2
$ws.Range("J129").Value = 1278355.8
$ws.Range("L129").Value = 3835067.4
$ws.Range("N129").Value = -3845067.4
$ws.Range("H138").Value = 3693.8284
$ws.Range("I138").Value = 1448.1482
$ws.Range("J138").Value = 4535.9585
$ws.Range("K138").Value = 4344.444600000001
$ws.Range("L138").Value = 13607.8755
$ws.Range("M138").Value = 795.5553999999993
$ws.Range("N138").Value = -23887.8755

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1712.0975
$ws.Range("I2").Value = 1649.0555
$ws.Range("J2").Value = 1761.4348
$ws.Range("K2").Value = 1649.0555
$ws.Range("L2").Value = 1761.4348
$ws.Range("M2").Value = -1536.0555
$ws.Range("N2").Value = -1987.4348
$ws.Range("H32").Value = 15162.603
$ws.Range("I32").Value = 11049.405
$ws.Range("K32").Value = 11049.405
$ws.Range("M32").Value = -10762.405
$ws.Range("H45").Value = 2272.7273
$ws.Range("I45").Value = 2960
$ws.Range("J45").Value = 1700
$ws.Range("K45").Value = 2960
$ws.Range("L45").Value = 1700
$ws.Range("M45").Value = -2583
$ws.Range("N45").Value = -2454
$ws.Range("H74").Value = 1147.3334
$ws.Range("I74").Value = 1137.1562
$ws.Range("J74").Value = 1179.9
$ws.Range("K74").Value = 1137.1562
$ws.Range("L74").Value = 1179.9
$ws.Range("M74").Value = -263.1561999999999
$ws.Range("N74").Value = -2927.9
$ws.Range("H77").Value = 1147.3334
$ws.Range("I77").Value = 1137.1562
$ws.Range("J77").Value = 1179.9
$ws.Range("K77").Value = 5685.780999999999
$ws.Range("L77").Value = 5899.5
$ws.Range("M77").Value = -1317.780999999999
$ws.Range("N77").Value = -14635.5
$ws.Range("H110").Value = 1329
$ws.Range("I110").Value = 1172.3334
$ws.Range("J110").Value = 1799
$ws.Range("K110").Value = 1172.3334
$ws.Range("L110").Value = 1799
$ws.Range("M110").Value = 872.6666
$ws.Range("N110").Value = -5889
$ws.Range("H116").Value = 1712.0975
$ws.Range("I116").Value = 1649.0555
$ws.Range("J116").Value = 1761.4348
$ws.Range("K116").Value = 1649.0555
$ws.Range("L116").Value = 1761.4348
$ws.Range("M116").Value = 644.9445000000001
$ws.Range("N116").Value = -6349.4348
$ws.Range("H122").Value = 1551.375
$ws.Range("I122").Value = 1568.5
$ws.Range("K122").Value = 4705.5
$ws.Range("M122").Value = -2255.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1712.0975
$ws.Range("I3").Value = 1649.0555
$ws.Range("J3").Value = 1761.4348
$ws.Range("K3").Value = 1649.0555
$ws.Range("L3").Value = 1761.4348
$ws.Range("M3").Value = -1535.0555
$ws.Range("N3").Value = -1989.4348
$ws.Range("H86").Value = 76935550
$ws.Range("I86").Value = 200022400
$ws.Range("J86").Value = 6273.375
$ws.Range("K86").Value = 200022400
$ws.Range("L86").Value = 6273.375
$ws.Range("M86").Value = -200021277
$ws.Range("N86").Value = -8519.375
$ws.Range("H89").Value = 76935550
$ws.Range("I89").Value = 200022400
$ws.Range("J89").Value = 6273.375
$ws.Range("K89").Value = 1000112000
$ws.Range("L89").Value = 31366.875
$ws.Range("M89").Value = -1000106384
$ws.Range("N89").Value = -42598.875
$ws.Range("H107").Value = 8602.352999999999
$ws.Range("I107").Value = 848
$ws.Range("J107").Value = 33804
$ws.Range("K107").Value = 848
$ws.Range("L107").Value = 33804
$ws.Range("M107").Value = 1072
$ws.Range("N107").Value = -37644

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 100006
$ws.Range("J11").Value = 100006
$ws.Range("L11").Value = 100006
$ws.Range("N11").Value = -100286
$ws.Range("H22").Value = 1004.5714
$ws.Range("I22").Value = 1251.3334
$ws.Range("J22").Value = 560.4
$ws.Range("K22").Value = 1251.3334
$ws.Range("L22").Value = 560.4
$ws.Range("M22").Value = -901.3334
$ws.Range("N22").Value = -1260.4
$ws.Range("H31").Value = 49607.5
$ws.Range("I31").Value = 4664.364
$ws.Range("J31").Value = 94550.63
$ws.Range("K31").Value = 4664.364
$ws.Range("L31").Value = 94550.63
$ws.Range("M31").Value = -4369.364
$ws.Range("N31").Value = -95140.63
$ws.Range("H34").Value = 49607.5
$ws.Range("I34").Value = 4664.364
$ws.Range("J34").Value = 94550.63
$ws.Range("K34").Value = 4664.364
$ws.Range("L34").Value = 94550.63
$ws.Range("M34").Value = -4462.364
$ws.Range("N34").Value = -94954.63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H122").Value = 1103.3334
$ws.Range("I122").Value = 1103.3334
$ws.Range("K122").Value = 3310.0002
$ws.Range("M122").Value = -860.0001999999999

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2020
$ws.Range("I62").Value = 1750
$ws.Range("J62").Value = 2200
$ws.Range("K62").Value = 5250
$ws.Range("L62").Value = 6600
$ws.Range("M62").Value = -4564
$ws.Range("N62").Value = -7972
$ws.Range("H65").Value = 2020
$ws.Range("I65").Value = 1750
$ws.Range("J65").Value = 2200
$ws.Range("K65").Value = 15750
$ws.Range("L65").Value = 19800
$ws.Range("M65").Value = -12318
$ws.Range("N65").Value = -26664
$ws.Range("H96").Value = 4998.2856
$ws.Range("J96").Value = 4998.2856
$ws.Range("L96").Value = 14994.8568
$ws.Range("N96").Value = -19112.8568
$ws.Range("H117").Value = 1859
$ws.Range("J117").Value = 2770.3
$ws.Range("L117").Value = 8310.900000000001
$ws.Range("N117").Value = -15194.9
$ws.Range("H121").Value = 834.5263
$ws.Range("I121").Value = 400
$ws.Range("J121").Value = 858.6667
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 2576.0001
$ws.Range("M121").Value = 110
$ws.Range("N121").Value = -5196.0001
$ws.Range("H123").Value = 1725.7273
$ws.Range("I123").Value = 1081.6666
$ws.Range("J123").Value = 2498.6
$ws.Range("K123").Value = 3244.9998
$ws.Range("L123").Value = 7495.799999999999
$ws.Range("M123").Value = -794.9998000000001
$ws.Range("N123").Value = -12395.8
$ws.Range("H131").Value = 12872942
$ws.Range("I131").Value = 45545890
$ws.Range("J131").Value = 37142.605
$ws.Range("K131").Value = 136637670
$ws.Range("L131").Value = 111427.815
$ws.Range("M131").Value = -136632630
$ws.Range("N131").Value = -121507.815
$ws.Range("H137").Value = 75734.86
$ws.Range("I137").Value = 3428.9
$ws.Range("J137").Value = 256499.75
$ws.Range("K137").Value = 10286.7
$ws.Range("L137").Value = 769499.25
$ws.Range("M137").Value = -5186.700000000001
$ws.Range("N137").Value = -779699.25

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 150
$ws.Range("I13").Value = 150
$ws.Range("K13").Value = 150
$ws.Range("M13").Value = -11
$ws.Range("H122").Value = 1645599.1
$ws.Range("I122").Value = 2193815.5
$ws.Range("J122").Value = 950
$ws.Range("K122").Value = 6581446.5
$ws.Range("L122").Value = 2850
$ws.Range("M122").Value = -6578996.5
$ws.Range("N122").Value = -7750
$ws.Range("H126").Value = 4437.231
$ws.Range("I126").Value = 3166.111
$ws.Range("K126").Value = 9498.332999999999
$ws.Range("M126").Value = -7028.332999999999
$ws.Range("H132").Value = 2680.3416
$ws.Range("I132").Value = 2274.3125
$ws.Range("J132").Value = 4124
$ws.Range("K132").Value = 6822.9375
$ws.Range("L132").Value = 12372
$ws.Range("M132").Value = -4292.9375
$ws.Range("N132").Value = -17432

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 28572142
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H61").Value = 1801.2759
$ws.Range("I61").Value = 1730.5
$ws.Range("J61").Value = 2141
$ws.Range("K61").Value = 1730.5
$ws.Range("L61").Value = 2141
$ws.Range("M61").Value = -1528.5
$ws.Range("N61").Value = -2545
$ws.Range("H113").Value = 1801.2759
$ws.Range("I113").Value = 1730.5
$ws.Range("J113").Value = 2141
$ws.Range("K113").Value = 1730.5
$ws.Range("L113").Value = 2141
$ws.Range("M113").Value = 439.5
$ws.Range("N113").Value = -6481
